# Refresh the "Cluster name" / "Active cases" table (Sheet1!A2:B..) with an
# updated data upload: some clusters were removed, several new clusters were
# added, a couple of names were corrected, and active-case counts changed
# throughout. The table grows from 90 data rows to 94 data rows (95 incl.
# header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 94,2
$arr[0,0] = '3535 Opal Meadow Heights Aged Care Community Meadow Heights'
$arr[0,1] = 28
$arr[1,0] = 'ABD Group 512 Melbourne Road Construction Site Spotswood'
$arr[1,1] = 5
$arr[2,0] = 'Acquire BPO Southbank'
$arr[2,1] = 7
$arr[3,0] = 'Adorn Cosmetics Clayton'
$arr[3,1] = 5
$arr[4,0] = 'Al Haj Halal Meats Glenroy'
$arr[4,1] = 34
$arr[5,0] = 'Al-Taqwa College Truganina'
$arr[5,1] = 5
$arr[6,0] = 'Best&Less Fountain Gate Narre Warren'
$arr[6,1] = 5
$arr[7,0] = 'Budget Car and Truck Rentals Campbellfield'
$arr[7,1] = 7
$arr[8,0] = 'CS Square Caroline Springs'
$arr[8,1] = 9
$arr[9,0] = 'Cannie Road Construction Site Cannie'
$arr[9,1] = 7
$arr[10,0] = 'Caroline Springs Police Station'
$arr[10,1] = 12
$arr[11,0] = 'Cedars Medical Clinic Coburg'
$arr[11,1] = 41
$arr[12,0] = 'Cedars Medical Clinic Coburg'
$arr[12,1] = 28
$arr[13,0] = 'Chemist Warehouse Campbellfield DC'
$arr[13,1] = 5
$arr[14,0] = 'Chemist Warehouse Fillo Drive Somerton'
$arr[14,1] = 5
$arr[15,0] = 'City of Moreland Community'
$arr[15,1] = 6
$arr[16,0] = 'City of Wyndham Community'
$arr[16,1] = 6
$arr[17,0] = 'Classy Cabinets and Kitchens Craigieburn'
$arr[17,1] = 10
$arr[18,0] = 'Coles Aurora Village Epping'
$arr[18,1] = 6
$arr[19,0] = 'Coles Broadmeadows Central Shopping Centre'
$arr[19,1] = 9
$arr[20,0] = 'Coles Campbellfield Plaza Campbellfield'
$arr[20,1] = 9
$arr[21,0] = 'Coles Coburg North Village'
$arr[21,1] = 26
$arr[22,0] = 'Coles Greenvale Shopping Centre'
$arr[22,1] = 7
$arr[23,0] = 'Coles Pakenham Place Shopping Centre'
$arr[23,1] = 11
$arr[24,0] = 'Coles Roxburgh Village Roxburgh Park'
$arr[24,1] = 11
$arr[25,0] = 'Community Kids Meadow Heights'
$arr[25,1] = 12
$arr[26,0] = 'Construction Site Olea Apartment Caulfield North'
$arr[26,1] = 16
$arr[27,0] = 'Costco Wholesale Epping'
$arr[27,1] = 24
$arr[28,0] = 'Crusader Caravans Epping'
$arr[28,1] = 14
$arr[29,0] = 'DRC Laverton Automotive Repairs Laverton North'
$arr[29,1] = 5
$arr[30,0] = 'Direct Freight Express Cambellfield'
$arr[30,1] = 13
$arr[31,0] = 'Epworth Healthcare Epworth Richmond Emergency Department'
$arr[31,1] = 5
$arr[32,0] = 'Fitzroy Community School Fitzroy North'
$arr[32,1] = 36
$arr[33,0] = 'Fonterra Manufacturing Workplace Campbellfield'
$arr[33,1] = 7
$arr[34,0] = 'Glenroy West Primary School'
$arr[34,1] = 6
$arr[35,0] = 'Goodstart Early Learning Altona'
$arr[35,1] = 5
$arr[36,0] = 'Green Leaves Early Learning Centre Highlands Craigieburn'
$arr[36,1] = 7
$arr[37,0] = 'Gumboots Early Learning Centre South Morang'
$arr[37,1] = 5
$arr[38,0] = 'Hamilton Marino 236 Jasper Road McKinnon'
$arr[38,1] = 12
$arr[39,0] = 'Health Care Providers Association South Melbourne'
$arr[39,1] = 13
$arr[40,0] = 'IGA Meadow Heights Shopping Centre Meadow Heights'
$arr[40,1] = 6
$arr[41,0] = 'ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine'
$arr[41,1] = 9
$arr[42,0] = 'Ibis Kingsgate Hotel Melbourne'
$arr[42,1] = 6
$arr[43,0] = 'Ilim College Glenroy Campus Hadfield'
$arr[43,1] = 16
$arr[44,0] = 'Ilim Learning Sanctuary Glenroy'
$arr[44,1] = 11
$arr[45,0] = 'Impact Designer Homes Epping'
$arr[45,1] = 5
$arr[46,0] = 'Industrial Galvanizers Valmont Coatings Campbellfield'
$arr[46,1] = 22
$arr[47,0] = 'Islamic College of Melbourne Tarneit'
$arr[47,1] = 5
$arr[48,0] = 'KFC Fawkner'
$arr[48,1] = 8
$arr[49,0] = 'Kasr Sweets Coolaroo'
$arr[49,1] = 5
$arr[50,0] = 'Kids House Early Learning Cheltenham'
$arr[50,1] = 12
$arr[51,0] = 'Learning Nest Early Learning Centre Meadow Heights'
$arr[51,1] = 5
$arr[52,0] = 'Level Crossing Removal Project Lilydale Construction Site John Street'
$arr[52,1] = 9
$arr[53,0] = 'Lineage Logistics Laverton North'
$arr[53,1] = 8
$arr[54,0] = 'Linfox Somerton National Distribution Centre Somerton'
$arr[54,1] = 9
$arr[55,0] = 'McDonalds Thomastown II'
$arr[55,1] = 7
$arr[56,0] = 'Melbourne Metropolitan Remand Centre Ravenhall'
$arr[56,1] = 11
$arr[57,0] = 'Melbourne Truck Repairs Campbellfield'
$arr[57,1] = 7
$arr[58,0] = 'Melbourne West Police Station Docklands'
$arr[58,1] = 7
$arr[59,0] = 'Melbourne Youth Justice Centre Parkville'
$arr[59,1] = 5
$arr[60,0] = 'Melton Police Station Melton'
$arr[60,1] = 5
$arr[61,0] = 'Mercy Hospital for Women Heidelberg'
$arr[61,1] = 5
$arr[62,0] = 'Mernda YMCA Early Learning Centre Mernda'
$arr[62,1] = 5
$arr[63,0] = 'Montessori Beginnings Greenvale'
$arr[63,1] = 5
$arr[64,0] = 'MyCentre Childcare Broadmeadows'
$arr[64,1] = 17
$arr[65,0] = 'National Gallery of Victoria Melbourne'
$arr[65,1] = 9
$arr[66,0] = 'Newbury Child and Community Centre Craigieburn'
$arr[66,1] = 5
$arr[67,0] = 'Nido Early School Moonee Ponds'
$arr[67,1] = 14
$arr[68,0] = 'North Geelong House Party'
$arr[68,1] = 7
$arr[69,0] = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$arr[69,1] = 47
$arr[70,0] = 'Northern Health The Northern Hospital Epping'
$arr[70,1] = 5
$arr[71,0] = 'OnQ Plumbing and Excavations Craigieburn'
$arr[71,1] = 18
$arr[72,0] = 'Oporto Coolaroo'
$arr[72,1] = 13
$arr[73,0] = 'Our Lady Help of Christian''s Primary School Brunswick East'
$arr[73,1] = 8
$arr[74,0] = 'Paisley Park Early Learning Centre Bundoora'
$arr[74,1] = 8
$arr[75,0] = 'Panorama Construction Site Whitehorse Rd Box Hill'
$arr[75,1] = 17
$arr[76,0] = 'Private Residence Northern Community Services Fawkner'
$arr[76,1] = 5
$arr[77,0] = 'Ramsay Health Care Warringal Private Hospital Heidelberg'
$arr[77,1] = 9
$arr[78,0] = 'Richmond Quarter 261-271 Bridge Road Construction Site Richmond'
$arr[78,1] = 10
$arr[79,0] = 'Salta Drive Construction Site Rangedale Drainage Altona North'
$arr[79,1] = 7
$arr[80,0] = 'Sharpline Stainless Steel Coburg North'
$arr[80,1] = 5
$arr[81,0] = 'St Vincents Hospital Emergency Department Melbourne'
$arr[81,1] = 5
$arr[82,0] = 'Tek Foods Somerton'
$arr[82,1] = 12
$arr[83,0] = 'The Homestead Child and Family Centre Roxburgh Park'
$arr[83,1] = 12
$arr[84,0] = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'
$arr[84,1] = 11
$arr[85,0] = 'ThorwestenCabinets Pakenham'
$arr[85,1] = 13
$arr[86,0] = 'Total Window Concepts Hoppers Crossing'
$arr[86,1] = 6
$arr[87,0] = 'Unilodge College Square Student Accommodation 570 Lygon Street Carlton'
$arr[87,1] = 14
$arr[88,0] = 'Werribee Mercy Hospital Emergency Department'
$arr[88,1] = 8
$arr[89,0] = 'Western Health Footscray Hospital Emergency Department'
$arr[89,1] = 6
$arr[90,0] = 'Western Health Sunshine Hospital Emergency Department'
$arr[90,1] = 8
$arr[91,0] = 'Woodlands Long Day Care and Kindergarten Roxburgh Park'
$arr[91,1] = 5
$arr[92,0] = 'Woolworths Greenvale Lakes Roxburgh Park'
$arr[92,1] = 5
$arr[93,0] = 'Yara Childcare Centre Truganina'
$arr[93,1] = 10

$ws.Range("A2:B95").Value = $arr
